$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Table" grows by one column (B2:E6 -> B2:F6), adding a "FINAL SUB" column.
$lo = $ws.ListObjects.Item(1)

# Preserve the border formatting that currently sits on F4 (the table's bottom-right
# corner marker, cell style index 1) by copying it onto what will become the new
# corner (G4) before the table expansion shifts things around.
[void]$ws.Range("F4").Copy()
[void]$ws.Range("G4").PasteSpecial(-4122)

# Add the new table column; Excel auto-names it and the table range grows to B2:F6.
$col = $lo.ListColumns.Add()

# Restore F4 to a plain (unbordered) format like its neighbours C4/D4/E4, since it is
# no longer the last column of the table.
[void]$ws.Range("C4").Copy()
[void]$ws.Range("F4").PasteSpecial(-4122)

# Header + data for the new "FINAL SUB" column.
$ws.Range("F2").Value = "FINAL SUB"
$ws.Range("F3").Value = 33
$ws.Range("F4").Value = 33
$ws.Range("F5").Value = 33
$ws.Range("F6").Value = 0

# Match the new column width used for the rest of the table columns.
$ws.Range("F1").ColumnWidth = 16.6667

# Move the active selection to F6, matching the saved cursor position.
[void]$ws.Range("F6").Select()
